# "Creando boleta de registro" - adds 5 new menu rows (36-40) to the Menu
# seed sheet: "Mis servicios" (+ its detail submenu) and "Autorización"
# (+ nuevo/editar submenus), mirroring the existing catalog-menu pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- numeric columns (A id, B parent id, H order) -------------------------
$ws.Cells.Item(36, 1).Value = 36
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 8).Value = 7

$ws.Cells.Item(37, 1).Value = 37
$ws.Cells.Item(37, 2).Value = 36
$ws.Cells.Item(37, 8).Value = 0

$ws.Cells.Item(38, 1).Value = 38
$ws.Cells.Item(38, 2).Value = 10
$ws.Cells.Item(38, 8).Value = 5

$ws.Cells.Item(39, 1).Value = 39
$ws.Cells.Item(39, 2).Value = 38
$ws.Cells.Item(39, 8).Value = 0

$ws.Cells.Item(40, 1).Value = 40
$ws.Cells.Item(40, 2).Value = 38
$ws.Cells.Item(40, 8).Value = 0

# --- text columns (C name, D slug, E icon, F state, G role, I description)
$ws.Range("C36").Value = "Mis servicios"
$ws.Range("D36").Value = "mis-servicios"
$ws.Range("G36").Value = "Regular"
$ws.Range("I36").Value = "Menú para visualizar servicios"

$ws.Range("C37").Value = "Mis servicios detalle"
$ws.Range("D37").Value = "mis-servicios/detalle"

$ws.Range("E36").Value = "assignment_ind"

$ws.Range("D38").Value = "autorizaciones"
$ws.Range("I38").Value = "Submenú para autorizaciones"
$ws.Range("C38").Value = "Autorización"

$ws.Range("C39").Value = "Autorización nuevo"
$ws.Range("D39").Value = "autorizaciones/nuevo"
$ws.Range("I39").Value = "Ruta para nueva autorización"

$ws.Range("C40").Value = "Autorización editar"
$ws.Range("D40").Value = "autorizaciones/editar"
$ws.Range("I40").Value = "Rutar para editar autorización"

$ws.Range("F36").Value = "visible"
$ws.Range("E37").Value = "minimize"
$ws.Range("F37").Value = "oculto"
$ws.Range("G37").Value = "Regular"
$ws.Range("I37").Value = "Ruta para detalle de servicio"

$ws.Range("E38").Value = "panorama_fish_eye"
$ws.Range("F38").Value = "visible"
$ws.Range("G38").Value = "Digitador"

$ws.Range("E39").Value = "minimize"
$ws.Range("F39").Value = "oculto"
$ws.Range("G39").Value = "Digitador"

$ws.Range("E40").Value = "minimize"
$ws.Range("F40").Value = "oculto"
$ws.Range("G40").Value = "Digitador"

# --- view: scroll/select like the author left it --------------------------
$ws.Range("E39").Select()
